$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 186.68182
$ws.Range("I33").Value = 89.0625
$ws.Range("J33").Value = 447
$ws.Range("K33").Value = 89.0625
$ws.Range("L33").Value = 447
$ws.Range("M33").Value = 139.9375
$ws.Range("N33").Value = -905
# Row 86
$ws.Range("H86").Value = 1418.9
$ws.Range("I86").Value = 1227
$ws.Range("J86").Value = 1866.6666
$ws.Range("K86").Value = 1227
$ws.Range("L86").Value = 1866.6666
$ws.Range("M86").Value = -104
$ws.Range("N86").Value = -4112.6666
# Row 89
$ws.Range("H89").Value = 1418.9
$ws.Range("I89").Value = 1227
$ws.Range("J89").Value = 1866.6666
$ws.Range("K89").Value = 6135
$ws.Range("L89").Value = 9333.333000000001
$ws.Range("M89").Value = -519
$ws.Range("N89").Value = -20565.333
# Row 98
$ws.Range("H98").Value = 350546.97
$ws.Range("I98").Value = 466312.1
$ws.Range("J98").Value = 3251.625
$ws.Range("K98").Value = 466312.1
$ws.Range("L98").Value = 3251.625
$ws.Range("M98").Value = -464814.1
$ws.Range("N98").Value = -6247.625
# Row 122
$ws.Range("H122").Value = 350546.97
$ws.Range("I122").Value = 466312.1
$ws.Range("J122").Value = 3251.625
$ws.Range("K122").Value = 1398936.3
$ws.Range("L122").Value = 9754.875
$ws.Range("M122").Value = -1396486.3
$ws.Range("N122").Value = -14654.875
# Row 135
$ws.Range("H135").Value = 2441.1177
$ws.Range("I135").Value = 2165.6428
$ws.Range("J135").Value = 3726.6667
$ws.Range("K135").Value = 19490.7852
$ws.Range("L135").Value = 33540.0003
$ws.Range("M135").Value = -16955.7852
$ws.Range("N135").Value = -38610.0003
# Row 138
$ws.Range("H138").Value = 6739662.5
$ws.Range("I138").Value = 1898563
$ws.Range("K138").Value = 5695689
$ws.Range("M138").Value = -5690549
# Row 141
$ws.Range("H141").Value = 1935.4166
$ws.Range("I141").Value = 1540.238
$ws.Range("K141").Value = 4620.714
$ws.Range("M141").Value = 559.2860000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 21184.51
$ws.Range("I32").Value = 1907.1063
$ws.Range("K32").Value = 1907.1063
$ws.Range("M32").Value = -1620.1063
# Row 45
$ws.Range("H45").Value = 905
$ws.Range("I45").Value = 722.6667
$ws.Range("J45").Value = 1178.5
$ws.Range("K45").Value = 722.6667
$ws.Range("L45").Value = 1178.5
$ws.Range("M45").Value = -345.6667
$ws.Range("N45").Value = -1932.5
# Row 61
$ws.Range("H61").Value = 2349.1177
$ws.Range("I61").Value = 1789.8928
$ws.Range("J61").Value = 4958.8335
$ws.Range("K61").Value = 1789.8928
$ws.Range("L61").Value = 4958.8335
$ws.Range("M61").Value = -1577.8928
$ws.Range("N61").Value = -5382.8335
# Row 88
$ws.Range("H88").Value = 6719.25
$ws.Range("J88").Value = 8192.333000000001
$ws.Range("L88").Value = 8192.333000000001
$ws.Range("N88").Value = -9004.333000000001
# Row 91
$ws.Range("H91").Value = 6719.25
$ws.Range("J91").Value = 8192.333000000001
$ws.Range("L91").Value = 8192.333000000001
$ws.Range("N91").Value = -11000.333
# Row 132
$ws.Range("H132").Value = 2407.7036
$ws.Range("I132").Value = 2088.4167
$ws.Range("K132").Value = 6265.250100000001
$ws.Range("M132").Value = -3735.250100000001
# Row 136
$ws.Range("H136").Value = 2349.1177
$ws.Range("I136").Value = 1789.8928
$ws.Range("J136").Value = 4958.8335
$ws.Range("K136").Value = 5369.678400000001
$ws.Range("L136").Value = 14876.5005
$ws.Range("M136").Value = -2819.678400000001
$ws.Range("N136").Value = -19976.5005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 8214.066000000001
$ws.Range("I86").Value = 1628.8572
$ws.Range("J86").Value = 13976.125
$ws.Range("K86").Value = 1628.8572
$ws.Range("L86").Value = 13976.125
$ws.Range("M86").Value = -505.8571999999999
$ws.Range("N86").Value = -16222.125
# Row 89
$ws.Range("H89").Value = 8214.066000000001
$ws.Range("I89").Value = 1628.8572
$ws.Range("J89").Value = 13976.125
$ws.Range("K89").Value = 8144.286
$ws.Range("L89").Value = 69880.625
$ws.Range("M89").Value = -2528.286
$ws.Range("N89").Value = -81112.625
# Row 134
$ws.Range("H134").Value = 5345.8667
$ws.Range("I134").Value = 3208.2856
$ws.Range("J134").Value = 7216.25
$ws.Range("K134").Value = 9624.856800000001
$ws.Range("L134").Value = 21648.75
$ws.Range("M134").Value = -7089.856800000001
$ws.Range("N134").Value = -26718.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1500.8214
$ws.Range("I31").Value = 1105.3914
$ws.Range("K31").Value = 1105.3914
$ws.Range("M31").Value = -810.3914
# Row 34
$ws.Range("H34").Value = 1500.8214
$ws.Range("I34").Value = 1105.3914
$ws.Range("K34").Value = 1105.3914
$ws.Range("M34").Value = -903.3914
# Row 107
$ws.Range("H107").Value = 390.2857
$ws.Range("I107").Value = 143.26666
$ws.Range("J107").Value = 675.3077
$ws.Range("K107").Value = 143.26666
$ws.Range("L107").Value = 675.3077
$ws.Range("M107").Value = 1776.73334
$ws.Range("N107").Value = -4515.3077
# Row 132
$ws.Range("H132").Value = 3079.0527
$ws.Range("I132").Value = 2326.5557
$ws.Range("J132").Value = 3756.3
$ws.Range("K132").Value = 6979.6671
$ws.Range("L132").Value = 11268.9
$ws.Range("M132").Value = -4449.6671
$ws.Range("N132").Value = -16328.9
# Row 134
$ws.Range("H134").Value = 3461.3333
$ws.Range("I134").Value = 1715.6
$ws.Range("J134").Value = 5048.364
$ws.Range("K134").Value = 5146.799999999999
$ws.Range("L134").Value = 15145.092
$ws.Range("M134").Value = -2611.799999999999
$ws.Range("N134").Value = -20215.092

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 24
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").ClearContents()
# Row 33
$ws.Range("H33").Value = 290.75
$ws.Range("I33").Value = 212.71428
$ws.Range("K33").Value = 1276.28568
$ws.Range("M33").Value = -993.28568
# Row 35
$ws.Range("H35").Value = 1500
$ws.Range("J35").Value = 1500
$ws.Range("L35").Value = 4500
$ws.Range("N35").Value = -5076
# Row 36
$ws.Range("H36").Value = 50500
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 50500
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 151500
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -151838
# Row 128
$ws.Range("H128").Value = 197805.25
$ws.Range("I128").Value = 197805.25
$ws.Range("K128").Value = 593415.75
$ws.Range("M128").Value = -588435.75
# Row 131
$ws.Range("H131").Value = 3174.9333
$ws.Range("I131").Value = 353.33334
$ws.Range("J131").Value = 3323.4385
$ws.Range("K131").Value = 1060.00002
$ws.Range("L131").Value = 9970.315500000001
$ws.Range("M131").Value = 3979.99998
$ws.Range("N131").Value = -20050.3155
# Row 140
$ws.Range("H140").Value = 15557.071
$ws.Range("I140").Value = 22811
$ws.Range("J140").Value = 2500
$ws.Range("K140").Value = 68433
$ws.Range("L140").Value = 7500
$ws.Range("M140").Value = -63253
$ws.Range("N140").Value = -17860

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5873.2856
$ws.Range("I70").Value = 5971.826
$ws.Range("J70").Value = 5420
$ws.Range("K70").Value = 5971.826
$ws.Range("L70").Value = 5420
$ws.Range("M70").Value = -5701.826
$ws.Range("N70").Value = -5960
# Row 73
$ws.Range("H73").Value = 5873.2856
$ws.Range("I73").Value = 5971.826
$ws.Range("J73").Value = 5420
$ws.Range("K73").Value = 5971.826
$ws.Range("L73").Value = 5420
$ws.Range("M73").Value = -5035.826
$ws.Range("N73").Value = -7292
# Row 122
$ws.Range("H122").Value = 1011030.94
$ws.Range("I122").Value = 1235491.5
$ws.Range("J122").Value = 958
$ws.Range("K122").Value = 3706474.5
$ws.Range("L122").Value = 2874
$ws.Range("M122").Value = -3704024.5
$ws.Range("N122").Value = -7774
# Row 132
$ws.Range("H132").Value = 2551.7407
$ws.Range("I132").Value = 2162.4583
$ws.Range("K132").Value = 6487.374899999999
$ws.Range("M132").Value = -3957.374899999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2903.4167
$ws.Range("I7").Value = 2018.2
$ws.Range("K7").Value = 2018.2
$ws.Range("M7").Value = -1906.2
# Row 126
$ws.Range("H126").Value = 2903.4167
$ws.Range("I126").Value = 2018.2
$ws.Range("K126").Value = 6054.6
$ws.Range("M126").Value = -3584.6
# Row 136
$ws.Range("H136").Value = 5250.276
$ws.Range("I136").Value = 2839.25
$ws.Range("J136").Value = 8217.691999999999
$ws.Range("K136").Value = 8517.75
$ws.Range("L136").Value = 24653.076
$ws.Range("M136").Value = -5967.75
$ws.Range("N136").Value = -29753.076

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 28704.666
$ws.Range("J64").Value = 29057
$ws.Range("L64").Value = 29057
$ws.Range("N64").Value = -29553
# Row 67
$ws.Range("H67").Value = 28704.666
$ws.Range("J67").Value = 29057
$ws.Range("L67").Value = 29057
$ws.Range("N67").Value = -30773
# Row 122
$ws.Range("H122").Value = 112628.22
$ws.Range("I122").Value = 251613.5
$ws.Range("J122").Value = 1440
$ws.Range("K122").Value = 754840.5
$ws.Range("L122").Value = 4320
$ws.Range("M122").Value = -752390.5
$ws.Range("N122").Value = -9220
# Row 136
$ws.Range("H136").Value = 20897296
$ws.Range("I136").Value = 23882180
$ws.Range("J136").Value = 3097.5
$ws.Range("K136").Value = 71646540
$ws.Range("L136").Value = 9292.5
$ws.Range("M136").Value = -71643990
$ws.Range("N136").Value = -14392.5
